$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row with two new columns (D1, E1), copying the
# formatting (bold font, border, centered alignment) from the
# existing header cells B1:C1.
$ws.Range("B1:C1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)
$ws.Range("D1").Value = 3
$ws.Range("E1").Value = 4

# Row 2
$ws.Range("C2").Value = -4.693393667675708
$ws.Range("D2").Value = -4.332345593343109
$ws.Range("E2").Value = -3.971195360520306

# Row 3
$ws.Range("C3").Value = -1.171375753587524
$ws.Range("D3").Value = -1.135929708662731
$ws.Range("E3").Value = -1.085622402784564

# Row 4
$ws.Range("C4").Value = -0.06850756529651124
$ws.Range("D4").Value = 0.004874729275821801
$ws.Range("E4").Value = 0.05974085484037735

# Row 5
$ws.Range("C5").Value = -0.4260526063017495
$ws.Range("D5").Value = -0.1731028680127843
$ws.Range("E5").Value = 0.02475535579086744

# Row 6
$ws.Range("C6").Value = -0.1257255649142397
$ws.Range("D6").Value = -0.1309219845841394
$ws.Range("E6").Value = -0.1313236579696304

# Row 7
$ws.Range("C7").Value = 0.01878386401140213
$ws.Range("D7").Value = 0.01141172071177196
$ws.Range("E7").Value = 0.00600914518709653
